# Update row 2 of the student list (ds_sinhvien) with the new student's
# data. The "cccd" id number is entered as text (it was already stored as
# text before the edit), so force a text format before assigning it,
# then restore the default "Normal" style so no new cell style is
# introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123456"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "Hoàng"
$ws.Range("D2").Value = "Nam"
$ws.Range("E2").Value = "nam_abc"
$ws.Range("F2").Value = "string"
$ws.Range("G2").Value = "string"
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 6
